# Update the two-digit multiplication problems in the table.
$d = $word.ActiveDocument

$pairs = @(
    @("39×53=", "14×68="),
    @("13×96=", "25×33="),
    @("34×95=", "28×94="),
    @("16×58=", "33×25="),
    @("64×55=", "40×27="),
    @("23×95=", "89×68="),
    @("48×71=", "90×27="),
    @("23×22=", "46×74="),
    @("66×38=", "33×80="),
    @("79×21=", "42×63="),
    @("75×83=", "12×37="),
    @("12×33=", "95×29="),
    @("78×76=", "74×32="),
    @("70×43=", "27×29="),
    @("72×63=", "27×34="),
    @("65×75=", "82×92="),
    @("59×73=", "87×85="),
    @("16×45=", "34×99="),
    @("30×65=", "74×89="),
    @("92×55=", "72×84="),
    @("63×58=", "81×59="),
    @("84×82=", "99×53="),
    @("41×64=", "56×57="),
    @("12×53=", "87×19="),
    @("57×94=", "90×54=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
